$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "imagen" column (E) values for rows 2-7 with the new image paths.
# Order matters for the shared-strings table build order, so assign in the
# same sequence the original edit was authored in.
$ws.Range("E2").Value = "/assets/img/productos/extintores/extintor-5-lbs-co2.png"
$ws.Range("E4").Value = "/assets/img/productos/extintores/Extintor-10-lbs-pqs.png"
$ws.Range("E3").Value = "/assets/img/productos/extintores/10-lbs-co2-1.png"
$ws.Range("E6").Value = "/assets/img/productos/extintores/EXTINTOR-5-LBS-PQS.png"
$ws.Range("E7").Value = "/assets/img/productos/extintores/EXTINTOR-20-LBS-PQS.png"
$ws.Range("E5").Value = "/assets/img/productos/extintores/1-11.png"

# Update the active selection to match the saved view state.
$ws.Range("F4").Select()
